$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row19 = New-Object 'object[,]' 1,24
$row19[0,0] = 4978
$row19[0,1] = 4801
$row19[0,2] = 4681
$row19[0,3] = 4625
$row19[0,4] = 4614
$row19[0,5] = 4647
$row19[0,6] = 4780
$row19[0,7] = 5025
$row19[0,8] = 5485
$row19[0,9] = 5857
$row19[0,10] = 6018
$row19[0,11] = 6107
$row19[0,12] = 6052
$row19[0,13] = 5915
$row19[0,14] = 5762
$row19[0,15] = 5650
$row19[0,16] = 5877
$row19[0,17] = 6104
$row19[0,18] = 6495
$row19[0,19] = 6502
$row19[0,20] = 6303
$row19[0,21] = 6064
$row19[0,22] = 5830
$row19[0,23] = 5464
$ws.Range("D19:AA19").Value = $row19

$row20 = New-Object 'object[,]' 1,24
$row20[0,0] = 5105
$row20[0,1] = 4906
$row20[0,2] = 4804
$row20[0,3] = 4817
$row20[0,4] = 4863
$row20[0,5] = 5035
$row20[0,6] = 5570
$row20[0,7] = 6085
$row20[0,8] = 6423
$row20[0,9] = 6506
$row20[0,10] = 6335
$row20[0,11] = 6265
$row20[0,12] = 6167
$row20[0,13] = 6090
$row20[0,14] = 5966
$row20[0,15] = 5907
$row20[0,16] = 5960
$row20[0,17] = 6123
$row20[0,18] = 6406
$row20[0,19] = 6470
$row20[0,20] = 6233
$row20[0,21] = 5931
$row20[0,22] = 5690
$row20[0,23] = 5329
$ws.Range("D20:AA20").Value = $row20

$row21 = New-Object 'object[,]' 1,24
$row21[0,0] = 4979
$row21[0,1] = 4737
$row21[0,2] = 4635
$row21[0,3] = 4623
$row21[0,4] = 4692
$row21[0,5] = 4812
$row21[0,6] = 5218
$row21[0,7] = 5674
$row21[0,8] = 6004
$row21[0,9] = 6050
$row21[0,10] = 5908
$row21[0,11] = 5774
$row21[0,12] = 5672
$row21[0,13] = 5612
$row21[0,14] = 5555
$row21[0,15] = 5548
$row21[0,16] = 5602
$row21[0,17] = 5750
$row21[0,18] = 6096
$row21[0,19] = 6151
$row21[0,20] = 5913
$row21[0,21] = 5636
$row21[0,22] = 5388
$row21[0,23] = 5015
$ws.Range("D21:AA21").Value = $row21

$row22 = New-Object 'object[,]' 1,24
$row22[0,0] = 4855
$row22[0,1] = 4614
$row22[0,2] = 4512
$row22[0,3] = 4501
$row22[0,4] = 4569
$row22[0,5] = 4689
$row22[0,6] = 5094
$row22[0,7] = 5549
$row22[0,8] = 5878
$row22[0,9] = 5924
$row22[0,10] = 5794
$row22[0,11] = 5671
$row22[0,12] = 5576
$row22[0,13] = 5521
$row22[0,14] = 5468
$row22[0,15] = 5462
$row22[0,16] = 5512
$row22[0,17] = 5649
$row22[0,18] = 5969
$row22[0,19] = 6020
$row22[0,20] = 5783
$row22[0,21] = 5506
$row22[0,22] = 5259
$row22[0,23] = 4886
$ws.Range("D22:AA22").Value = $row22

$row23 = New-Object 'object[,]' 1,24
$row23[0,0] = 4897
$row23[0,1] = 4656
$row23[0,2] = 4554
$row23[0,3] = 4542
$row23[0,4] = 4611
$row23[0,5] = 4730
$row23[0,6] = 5136
$row23[0,7] = 5591
$row23[0,8] = 5921
$row23[0,9] = 5967
$row23[0,10] = 5832
$row23[0,11] = 5706
$row23[0,12] = 5608
$row23[0,13] = 5552
$row23[0,14] = 5498
$row23[0,15] = 5491
$row23[0,16] = 5543
$row23[0,17] = 5683
$row23[0,18] = 6012
$row23[0,19] = 6065
$row23[0,20] = 5827
$row23[0,21] = 5550
$row23[0,22] = 5302
$row23[0,23] = 4930
$ws.Range("D23:AA23").Value = $row23

$row24 = New-Object 'object[,]' 1,27
$row24[0,0] = 2026
$row24[0,1] = 2
$row24[0,2] = 27
$row24[0,3] = 4930
$row24[0,4] = 4689
$row24[0,5] = 4587
$row24[0,6] = 4575
$row24[0,7] = 4643
$row24[0,8] = 4763
$row24[0,9] = 5169
$row24[0,10] = 5625
$row24[0,11] = 5954
$row24[0,12] = 6000
$row24[0,13] = 5863
$row24[0,14] = 5734
$row24[0,15] = 5634
$row24[0,16] = 5576
$row24[0,17] = 5521
$row24[0,18] = 5514
$row24[0,19] = 5560
$row24[0,20] = 5683
$row24[0,21] = 5973
$row24[0,22] = 6020
$row24[0,23] = 5789
$row24[0,24] = 5519
$row24[0,25] = 5279
$row24[0,26] = 4916
$ws.Range("A24:AA24").Value = $row24
